$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 17:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1388755
$ws.Range("C4").Value = 2921
$ws.Range("D4").Value = 262328
$ws.Range("E4").Value = 1044408
$ws.Range("G4").Value = 224
$ws.Range("H4").Value = 82019

# Row 7 - Reino Unido
$ws.Range("B7").Value = 226463
$ws.Range("C7").Value = 3403
$ws.Range("E7").Value = 193427
$ws.Range("G7").Value = 627
$ws.Range("H7").Value = 32692

# Row 19 - Paises Bajos
$ws.Range("F19").Value = 463

# Row 23 previously held Suiza; it now becomes Chile with the updated figures
$ws.Range("A23").Value = "Chile"
$ws.Range("B23").Value = 31721
$ws.Range("C23").Value = 1658
$ws.Range("D23").Value = 14125
$ws.Range("E23").Value = 17261
$ws.Range("F23").Value = 494
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = 335

# Row 24 previously held Chile; it now becomes Suiza, keeping the old Suiza figures
$ws.Range("A24").Value = "Suiza"
$ws.Range("B24").Value = 30380
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 26800
$ws.Range("E24").Value = 1723
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 12
$ws.Range("H24").Value = 1857

# Row 80 - Bulgaria
$ws.Range("B80").Value = 2023
$ws.Range("C80").Value = 33
$ws.Range("E80").Value = 1452
$ws.Range("F80").Value = 51
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 95
